# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46073
$ws.Range("B2").Value = 24.94
$ws.Range("C2").Value = 20.38
$ws.Range("D2").Value = 12.45
$ws.Range("E2").Value = 7.26
$ws.Range("F2").Value = 6.89
$ws.Range("G2").Value = 9.800000000000001
$ws.Range("H2").Value = 24.03
$ws.Range("I2").Value = 43.2
$ws.Range("J2").Value = 54.29
$ws.Range("K2").Value = 11.74
$ws.Range("L2").Value = 7.02
$ws.Range("M2").Value = 6.89
$ws.Range("N2").Value = 6.86
$ws.Range("O2").Value = 6.66
$ws.Range("P2").Value = 6.53
$ws.Range("Q2").Value = 6.7
$ws.Range("R2").Value = 6.86
$ws.Range("S2").Value = 12.48
$ws.Range("T2").Value = 59.1
$ws.Range("U2").Value = 93.45999999999999
$ws.Range("V2").Value = 100.54
$ws.Range("W2").Value = 96.02
$ws.Range("X2").Value = 60.07
$ws.Range("Y2").Value = 38.68
$ws.Range("Z2").Value = 30.12
$ws.Range("AB2").Value = 73.83
$ws.Range("AD2").Value = 98.28
$ws.Range("AF2").Value = 76.28
